$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 16.83279933333333
$ws.Range("H2").Value = 50.498398
$ws.Range("I2").Value = 0.04383102208811961
$ws.Range("J2").Value = 0.04383102208811961
$ws.Range("M2").Value = 145.7007446666667
$ws.Range("N2").Value = 437.1022340000001
$ws.Range("O2").Value = 0.2865937750105843
$ws.Range("P2").Value = 0.2865937750105843
$ws.Range("Q2").Value = 2452.551397691237
$ws.Range("R2").Value = 22072.96257922114
$ws.Range("S2").Value = 0.0125616980828065
$ws.Range("T2").Value = 0.0125616980828065

$ws.Range("G3").Value = 16.83279933333333
$ws.Range("H3").Value = 50.498398
$ws.Range("I3").Value = 0.04383102208811961
$ws.Range("J3").Value = 0.04383102208811961
$ws.Range("M3").Value = 168.7997026666667
$ws.Range("N3").Value = 506.3991080000001
$ws.Range("O3").Value = 0.3320294904365841
$ws.Range("P3").Value = 0.3320294904365841
$ws.Range("Q3").Value = 2841.371522514332
$ws.Range("R3").Value = 25572.34370262899
$ws.Range("S3").Value = 0.01455319192923302
$ws.Range("T3").Value = 0.01455319192923302

$ws.Range("G4").Value = 16.83279933333333
$ws.Range("H4").Value = 50.498398
$ws.Range("I4").Value = 0.04383102208811961
$ws.Range("J4").Value = 0.04383102208811961
$ws.Range("M4").Value = 128.1261546666667
$ws.Range("N4").Value = 384.378464
$ws.Range("O4").Value = 0.2520245069956105
$ws.Range("P4").Value = 0.2520245069956105
$ws.Range("Q4").Value = 2156.721850855631
$ws.Range("R4").Value = 19410.49665770067
$ws.Range("S4").Value = 0.01104649173287206
$ws.Range("T4").Value = 0.01104649173287206

$ws.Range("G5").Value = 16.83279933333333
$ws.Range("H5").Value = 50.498398
$ws.Range("I5").Value = 0.04383102208811961
$ws.Range("J5").Value = 0.04383102208811961
$ws.Range("M5").Value = 65.761079
$ws.Range("N5").Value = 197.283237
$ws.Range("O5").Value = 0.1293522275572212
$ws.Range("P5").Value = 0.1293522275572212
$ws.Range("Q5").Value = 1106.943046750481
$ws.Range("R5").Value = 9962.487420754325
$ws.Range("S5").Value = 0.005669640343208038
$ws.Range("T5").Value = 0.005669640343208038

$ws.Range("G6").Value = 332.3726806666667
$ws.Range("H6").Value = 997.1180420000001
$ws.Range("I6").Value = 0.8654671168650652
$ws.Range("J6").Value = 0.8654671168650654
$ws.Range("M6").Value = 145.7007446666667
$ws.Range("N6").Value = 437.1022340000001
$ws.Range("O6").Value = 0.2865937750105843
$ws.Range("P6").Value = 0.2865937750105843
$ws.Range("Q6").Value = 48426.94707998954
$ws.Range("R6").Value = 435842.5237199059
$ws.Range("S6").Value = 0.2480374881698855
$ws.Range("T6").Value = 0.2480374881698856

$ws.Range("G7").Value = 332.3726806666667
$ws.Range("H7").Value = 997.1180420000001
$ws.Range("I7").Value = 0.8654671168650652
$ws.Range("J7").Value = 0.8654671168650654
$ws.Range("M7").Value = 168.7997026666667
$ws.Range("N7").Value = 506.3991080000001
$ws.Range("O7").Value = 0.3320294904365841
$ws.Range("P7").Value = 0.3320294904365841
$ws.Range("Q7").Value = 56104.40967105629
$ws.Range("R7").Value = 504939.6870395066
$ws.Range("S7").Value = 0.2873606058023272
$ws.Range("T7").Value = 0.2873606058023273

$ws.Range("G8").Value = 332.3726806666667
$ws.Range("H8").Value = 997.1180420000001
$ws.Range("I8").Value = 0.8654671168650652
$ws.Range("J8").Value = 0.8654671168650654
$ws.Range("M8").Value = 128.1261546666667
$ws.Range("N8").Value = 384.378464
$ws.Range("O8").Value = 0.2520245069956105
$ws.Range("P8").Value = 0.2520245069956105
$ws.Range("Q8").Value = 42585.63349007194
$ws.Range("R8").Value = 383270.7014106475
$ws.Range("S8").Value = 0.2181189234488304
$ws.Range("T8").Value = 0.2181189234488305

$ws.Range("G9").Value = 332.3726806666667
$ws.Range("H9").Value = 997.1180420000001
$ws.Range("I9").Value = 0.8654671168650652
$ws.Range("J9").Value = 0.8654671168650654
$ws.Range("M9").Value = 65.761079
$ws.Range("N9").Value = 197.283237
$ws.Range("O9").Value = 0.1293522275572212
$ws.Range("P9").Value = 0.1293522275572212
$ws.Range("Q9").Value = 21857.18611076244
$ws.Range("R9").Value = 196714.674996862
$ws.Range("S9").Value = 0.1119500994440221
$ws.Range("T9").Value = 0.1119500994440221

$ws.Range("G10").Value = 34.50825133333333
$ws.Range("H10").Value = 103.524754
$ws.Range("I10").Value = 0.08985623225594501
$ws.Range("J10").Value = 0.08985623225594502
$ws.Range("M10").Value = 145.7007446666667
$ws.Range("N10").Value = 437.1022340000001
$ws.Range("O10").Value = 0.2865937750105843
$ws.Range("P10").Value = 0.2865937750105843
$ws.Range("Q10").Value = 5027.87791641116
$ws.Range("R10").Value = 45250.90124770044
$ws.Range("S10").Value = 0.02575223681045911
$ws.Range("T10").Value = 0.02575223681045911

$ws.Range("G11").Value = 34.50825133333333
$ws.Range("H11").Value = 103.524754
$ws.Range("I11").Value = 0.08985623225594501
$ws.Range("J11").Value = 0.08985623225594502
$ws.Range("M11").Value = 168.7997026666667
$ws.Range("N11").Value = 506.3991080000001
$ws.Range("O11").Value = 0.3320294904365841
$ws.Range("P11").Value = 0.3320294904365841
$ws.Range("Q11").Value = 5824.982564613271
$ws.Range("R11").Value = 52424.84308151944
$ws.Range("S11").Value = 0.02983491900849278
$ws.Range("T11").Value = 0.02983491900849278

$ws.Range("G12").Value = 34.50825133333333
$ws.Range("H12").Value = 103.524754
$ws.Range("I12").Value = 0.08985623225594501
$ws.Range("J12").Value = 0.08985623225594502
$ws.Range("M12").Value = 128.1261546666667
$ws.Range("N12").Value = 384.378464
$ws.Range("O12").Value = 0.2520245069956105
$ws.Range("P12").Value = 0.2520245069956105
$ws.Range("Q12").Value = 4421.409547610873
$ws.Range("R12").Value = 39792.68592849786
$ws.Range("S12").Value = 0.02264597263478761
$ws.Range("T12").Value = 0.02264597263478762

$ws.Range("G13").Value = 34.50825133333333
$ws.Range("H13").Value = 103.524754
$ws.Range("I13").Value = 0.08985623225594501
$ws.Range("J13").Value = 0.08985623225594502
$ws.Range("M13").Value = 65.761079
$ws.Range("N13").Value = 197.283237
$ws.Range("O13").Value = 0.1293522275572212
$ws.Range("P13").Value = 0.1293522275572212
$ws.Range("Q13").Value = 2269.299842083189
$ws.Range("R13").Value = 20423.6985787487
$ws.Range("S13").Value = 0.01162310380220552
$ws.Range("T13").Value = 0.01162310380220552

$ws.Range("G14").Value = 0.324754
$ws.Range("H14").Value = 0.974262
$ws.Range("I14").Value = 0.0008456287908700705
$ws.Range("J14").Value = 0.0008456287908700706
$ws.Range("M14").Value = 145.7007446666667
$ws.Range("N14").Value = 437.1022340000001
$ws.Range("O14").Value = 0.2865937750105843
$ws.Range("P14").Value = 0.2865937750105843
$ws.Range("Q14").Value = 47.31689963347867
$ws.Range("R14").Value = 425.852096701308
$ws.Range("S14").Value = 0.0002423519474330894
$ws.Range("T14").Value = 0.0002423519474330894

$ws.Range("G15").Value = 0.324754
$ws.Range("H15").Value = 0.974262
$ws.Range("I15").Value = 0.0008456287908700705
$ws.Range("J15").Value = 0.0008456287908700706
$ws.Range("M15").Value = 168.7997026666667
$ws.Range("N15").Value = 506.3991080000001
$ws.Range("O15").Value = 0.3320294904365841
$ws.Range("P15").Value = 0.3320294904365841
$ws.Range("Q15").Value = 54.81837863981067
$ws.Range("R15").Value = 493.365407758296
$ws.Range("S15").Value = 0.0002807736965310943
$ws.Range("T15").Value = 0.0002807736965310943

$ws.Range("G16").Value = 0.324754
$ws.Range("H16").Value = 0.974262
$ws.Range("I16").Value = 0.0008456287908700705
$ws.Range("J16").Value = 0.0008456287908700706
$ws.Range("M16").Value = 128.1261546666667
$ws.Range("N16").Value = 384.378464
$ws.Range("O16").Value = 0.2520245069956105
$ws.Range("P16").Value = 0.2520245069956105
$ws.Range("Q16").Value = 41.60948123261867
$ws.Range("R16").Value = 374.485331093568
$ws.Range("S16").Value = 0.0002131191791203237
$ws.Range("T16").Value = 0.0002131191791203237

$ws.Range("G17").Value = 0.324754
$ws.Range("H17").Value = 0.974262
$ws.Range("I17").Value = 0.0008456287908700705
$ws.Range("J17").Value = 0.0008456287908700706
$ws.Range("M17").Value = 65.761079
$ws.Range("N17").Value = 197.283237
$ws.Range("O17").Value = 0.1293522275572212
$ws.Range("P17").Value = 0.1293522275572212
$ws.Range("Q17").Value = 21.356173449566
$ws.Range("R17").Value = 192.205561046094
$ws.Range("S17").Value = 0.0001093839677855632
$ws.Range("T17").Value = 0.0001093839677855632

